$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTECreation")
Write-Output $ws.Name
